$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I1").Value = "pixel_size_mm"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I2").Value = 1.818
$ws.Range("F16").Select() | Out-Null
